$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 3 with the new test-log entry
$ws.Range("A3").Value = "Ball hitting corner of bricks"
$ws.Range("B3").Value = "ball going past corners of bricks breaking them"
$ws.Range("C3").Value = "When ball hits corner it should break the brick and its x or y speed should change depending on what direction it is going. Ef if going to bottome right should bounce to bottome left of canvas"
$ws.Range("D3").Value = "Sometimes the ball skims the corner breaking the brick but continuign in a stratight line where it came from"
$ws.Range("E3").Value = "The prolem was that it thought that it was coliding with other side of brick. To fix this added and if line to each collision checking what direcion the ball is going so it always bounces if it hits."

# Match the row height used in the source workbook for this row
$ws.Rows.Item(3).RowHeight = 66

# Update the selected cell to match the author's saved selection
$ws.Range("D6").Select()

$wb.Save()
